$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (new values) for rows 69-82, columns A,B,D,E,F,G,H,Q,R
# Derived by tracing the permutation of rows described in the diff.

$data = @{
    69 = @{ A=111785202; B=78512;  D="LC"; E=6456;   F="Skinnlav";        G="Leptogium saturninum";            H="(Dicks.) Nyl.";                  Q=577215.0430418774; R=6944631.445974576 }
    70 = @{ A=111785201; B=78512;  D="LC"; E=6456;   F="Skinnlav";        G="Leptogium saturninum";            H="(Dicks.) Nyl.";                  Q=577248.2772659193; R=6944530.940753835 }
    71 = @{ A=111785251; B=93161;  D="VU"; E=1079;   F="Aspfjädermossa";  G="Neckera pennata";                 H="Hedw.";                          Q=577283.2535308318; R=6944533.598891968 }
    72 = @{ A=111785229; B=78578;  D="NT"; E=6458;   F="Lunglav";         G="Lobaria pulmonaria";              H="(L.) Hoffm.";                    Q=577208.3826684169; R=6944521.722980071 }
    73 = @{ A=111785206; B=77268;  D="NT"; E=228912; F="Mörk kolflarnlav"; G="Carbonicola myrmecina";          H="(Ach.) Bendiksby & Timdal";      Q=577235.6798241453; R=6944655.86623876  }
    74 = @{ A=111785230; B=78578;  D="NT"; E=6458;   F="Lunglav";         G="Lobaria pulmonaria";              H="(L.) Hoffm.";                    Q=577261.8704127767; R=6944620.109213427 }
    75 = @{ A=111785199; B=89416;  D="LC"; E=1205;   F="Stor aspticka";   G="Phellinus populicola";            H="Niemelä";                         Q=577256.110519147;  R=6944531.123615563 }
    76 = @{ A=111785192; B=89405;  D="NT"; E=1202;   F="Ullticka";        G="Phellinidium ferrugineofuscum";   H="(P.Karst.) Fiasson & Niemelä";  Q=577281.7951240344; R=6944714.487089146 }
    77 = @{ A=111785228; B=78578;  D="NT"; E=6458;   F="Lunglav";         G="Lobaria pulmonaria";              H="(L.) Hoffm.";                    Q=577256.110519147;  R=6944531.123615563 }
    78 = @{ A=111785191; B=89405;  D="NT"; E=1202;   F="Ullticka";        G="Phellinidium ferrugineofuscum";   H="(P.Karst.) Fiasson & Niemelä";  Q=577235.6798241453; R=6944655.86623876  }
    79 = @{ A=111785244; B=96348;  D="VU"; E=220787; F="Knärot";          G="Goodyera repens";                 H="(L.) R. Br.";                    Q=577364.1822193989; R=6944621.578847388 }
    80 = @{ A=111785190; B=94134;  D="NT"; E=53;     F="Vedtrappmossa";   G="Crossocalyx hellerianus";         H="(Nees ex Lindenb.) Meyl.";       Q=577242.972410051;  R=6944540.945152098 }
    81 = @{ A=111785235; B=77267;  D="NT"; E=6446;   F="Kolflarnlav";     G="Carbonicola anthracophila";       H="(Nyl.) Bendiksby & Timdal";      Q=577226.625646919;  R=6944648.749557905 }
    82 = @{ A=111785200; B=78512;  D="LC"; E=6456;   F="Skinnlav";        G="Leptogium saturninum";            H="(Dicks.) Nyl.";                  Q=577256.110519147;  R=6944531.123615563 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
}
